$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force "Text" number format on cells whose new values look numeric, so Excel
# keeps them stored as text (matching the original inlineStr string cells)
# instead of silently converting them to numbers (which would also drop
# meaningful trailing zeros, e.g. "236.00" -> 236).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.880.95"
$ws.Range("E2").Value = "  +3.71%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.259.85"
$ws.Range("E3").Value = "  +3.24%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "253.43"
$ws.Range("E5").Value = "  -0.78%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.636"
$ws.Range("E6").Value = "  +1.26%  "

# Row 7 - Solana
$ws.Range("D7").Value = "71.59"
$ws.Range("E7").Value = "  +4.81%  "

# Row 8 - USDC
$ws.Range("D8").Value = "0.999"

# Row 9 - Cardano
$ws.Range("E9").Value = "  +11.81%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "41.12"
$ws.Range("E10").Value = "  +8.76%  "

# Row 11 - OKB
$ws.Range("D11").Value = "59.77"
$ws.Range("E11").Value = "  +1.13%  "

# Row 12 - Dogecoin
$ws.Range("D12").Value = "0.0968"
$ws.Range("E12").Value = "  +3.35%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "7.35"
$ws.Range("E13").Value = "  +2.94%  "

# Row 14 - TRON
$ws.Range("D14").Value = "0.105"
$ws.Range("E14").Value = "  +0.69%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.596.12"
$ws.Range("E15").Value = "  +3.46%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "0.888"
$ws.Range("E16").Value = "  +1.56%  "

# Row 17 - Chainlink
$ws.Range("D17").Value = "14.84"
$ws.Range("E17").Value = "  +2.48%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.250.14"
$ws.Range("E18").Value = "  +4.68%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "42.808.31"
$ws.Range("E19").Value = "  +3.73%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0981"
$ws.Range("E20").Value = "  +2.71%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +1.36%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "73.18"
$ws.Range("E22").Value = "  +1.76%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "236.00"
$ws.Range("E23").Value = "  +1.48%  "

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  +4.27%  "

# Row 25 - WEMIXToken
$ws.Range("E25").Value = "  +0.97%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "11.74"
$ws.Range("E26").Value = "  +0.04%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  -0.06%  "

# Row 28 - PancakeSwap
$ws.Range("E28").Value = "  -2.63%  "

# Row 29 - LEO
$ws.Range("E29").Value = "  -2.25%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  -0.90%  "

# Row 31 - Monero
$ws.Range("D31").Value = "167.87"
$ws.Range("E31").Value = "  -0.54%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "21.03"
$ws.Range("E32").Value = "  +1.70%  "

# Row 33 - Kaspa
$ws.Range("D33").Value = "0.131"
$ws.Range("E33").Value = "  +11.88%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").Value = "6.15"
$ws.Range("E34").Value = "  +12.36%  "

# Row 35 - Hedera
$ws.Range("D35").Value = "0.0790"
$ws.Range("E35").Value = "  +4.45%  "

# Row 36 - Stellar
$ws.Range("E36").Value = "  +1.32%  "

# Row 37 - InjectiveProtocol
$ws.Range("D37").Value = "29.33"
$ws.Range("E37").Value = "  +11.44%  "

# Row 38 - Filecoin
$ws.Range("E38").Value = "  +1.81%  "

# Row 39 - RenderToken
$ws.Range("D39").Value = "4.16"
$ws.Range("E39").Value = "  +0.31%  "

# Row 40 - VeChain
$ws.Range("D40").Value = "0.0320"
$ws.Range("E40").Value = "  +6.71%  "

# Row 41 - was LidoDAOToken, now THORChain
$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").Value = "6.22"
$ws.Range("E41").Value = "  +9.71%  "

# Row 42 - was THORChain, now LidoDAOToken
$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D42").Value = "2.30"
$ws.Range("E42").Value = "  +4.48%  "

# Row 43 - Celestia
$ws.Range("D43").Value = "12.51"
$ws.Range("E43").Value = "  -0.90%  "

# Row 44 - MultiversX
$ws.Range("D44").Value = "64.78"

# Row 45 - FTXToken
$ws.Range("D45").Value = "5.00"
$ws.Range("E45").Value = "  -2.79%  "

# Row 46 - Algorand
$ws.Range("E46").Value = "  +0.18%  "

# Row 47 - FraxShare
$ws.Range("D47").Value = "8.98"
$ws.Range("E47").Value = "  +3.90%  "

# Row 48 - Cronos
$ws.Range("E48").Value = "  +1.17%  "

# Row 49 - ARBITRUM
$ws.Range("E49").Value = "  +4.59%  "

# Row 50 - BinanceUSD
$ws.Range("E50").Value = "  -0.37%  "

# Row 51 - SynthetixNetwork
$ws.Range("E51").Value = "  +3.90%  "
